# "Added last minute updates"
# - Add a paragraph border (space-only, no visible line) to the first
#   paragraph (the hidden **ID__...__ID** marker paragraph).
# - Bump that paragraph's left indent from 120 -> 225 twips (11.25pt).
# - Update the marker text from the "topic_18" placeholder to the
#   "407_3" placeholder, and drop the now-redundant trailing-space run
#   that used to follow it.

$d = $word.ActiveDocument

# The marker paragraph is the very first paragraph in the document body.
$p = $d.Paragraphs(1)

# Add w:pBdr (top/left/bottom/right, each w:space="5", no line drawn) to
# this paragraph's pPr.
$borders = $p.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# w:ind w:left="120" -> w:ind w:left="225" (twips = points * 20).
$p.Range.ParagraphFormat.LeftIndent = 11.25

# Replace "**ID__AFFARS_5315_topic_18__ID** " (marker run + trailing
# space run) with "**ID__AFFARS_5315_407_3__ID**" (single run, no
# trailing space) in one shot so the orphan space run is removed rather
# than merely left behind.
$d.Content.Find.Execute("**ID__AFFARS_5315_topic_18__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_5315_407_3__ID**", 2)
